$d = $word.ActiveDocument

# Locate the paragraph that currently holds "Guardar productos favoritos en la cuenta".
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Guardar productos favoritos en la cuenta*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not find the 'Guardar productos favoritos en la cuenta' paragraph."
}

# Remove the two empty paragraphs immediately preceding it (one at a time -
# deleting a single paragraph's range is clean; the six empty paragraphs
# further up stay untouched). Each delete shifts later indices down by one.
$d.Paragraphs.Item($targetIndex - 1).Range.Delete()
$targetIndex = $targetIndex - 1
$d.Paragraphs.Item($targetIndex - 1).Range.Delete()
$targetIndex = $targetIndex - 1

# Replace the (now lone) "Guardar productos..." paragraph with the new
# business-description paragraphs, each separated by a blank paragraph.
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>La Barbería Golden, ubicada en el centro de Quilpué, es un establecimiento que ofrece una amplia gama de servicios estéticos y de cuidado personal, principalmente orientados al público masculino. Con un equipo reducido de 4 trabajadores, incluyendo al dueño que también es barbero, la barbería opera bajo un sistema de comisiones y permite a sus empleados manejar sus propios horarios.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>La base de clientes de la Barbería Golden se divide en un 20% de clientes de paso y un 80% de clientes fieles al local y/o a los barberos que trabajan en él. De los clientes fieles, un 70% son leales a la barbería en general, mientras que un 30% tienen una relación más estrecha con un barbero en particular.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Para mejorar su sistema actual de agendado de citas, la barbería tiene la intención de implementar un sistema web que permita a los barberos gestionar su propia disponibilidad a través de una página en línea. Esto facilitaría el proceso de reserva de citas tanto para los barberos como para los clientes.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>La página web también desempeñará un papel importante al explicar detalladamente los procesos de cada servicio ofrecido. Esto ayudará a los clientes a comprender mejor cómo funcionan los diferentes servicios y, en última instancia, a tomar decisiones informadas sobre cuál servicio es el más adecuado para sus necesidades.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Además de los servicios de barbería estándar, la Barbería Golden vende productos estéticos y consumibles en su local. Esta faceta de la barbería es de gran importancia para ellos, ya que contribuye a mejorar la experiencia y comodidad del cliente. La página web podría destacar estos productos y ofrecer la posibilidad de realizar compras en línea o proporcionar información detallada sobre los productos disponibles en la tienda física.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>En momentos de festividades u ocasiones especiales, la barbería también planea aplicar un aumento en el precio de los servicios, ya que se espera un aumento en el flujo de clientes. Este detalle puede ser comunicado claramente en la página web, junto con cualquier otra promoción o descuento que la barbería desee ofrecer.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>En resumen, la Barbería Golden busca modernizar su operación mediante la implementación de un sistema web de agendado de citas, la explicación detallada de sus servicios en línea y la promoción de sus productos estéticos. Esto contribuirá a brindar una experiencia de cliente más completa y conveniente, al tiempo que mejorará la eficiencia operativa de la barbería.</w:t></w:r></w:p>
'@

$d.Paragraphs.Item($targetIndex).Range.InsertXML($newXml)

# The run that hosts the picture gets marked NoProof (Word does this
# automatically for inline drawings) -> adds <w:rPr><w:noProof/></w:rPr>.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.NoProofing = 1
    }
}

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
